# Update workbook to reflect data through 2022-07-21 (commit: "Add data for 2022-07-29")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-07-21"

# Update the row label for July in column A (row 8)
$ws.Range("A8").Value = "July (through 07-21)"

# Update the July row (row 8) values for columns C..I
$ws.Range("C8").Value = 42
$ws.Range("D8").Value = 44
$ws.Range("E8").Value = 53
$ws.Range("F8").Value = 34
$ws.Range("G8").Value = 88
$ws.Range("H8").Value = 103
$ws.Range("I8").Value = 118

# Update the Total row (row 9) values for columns C..I
$ws.Range("C9").Value = 290
$ws.Range("D9").Value = 434
$ws.Range("E9").Value = 406
$ws.Range("F9").Value = 285
$ws.Range("G9").Value = 560
$ws.Range("H9").Value = 863
$ws.Range("I9").Value = 924
